$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.941.94'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '2.266.33'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.656'
$ws.Range("E5").Value = '  +4.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '233.32'
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.79'
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  +3.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0979'
$ws.Range("E10").Value = '  -7.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.03'
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.59'
$ws.Range("E12").Value = '  +2.54%  '
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").Value = '2.599.45'
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.14'
$ws.Range("E16").Value = '  +3.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  +1.91%  '
$ws.Range("D18").Value = '2.263.97'
$ws.Range("E18").Value = '  -0.41%  '
$ws.Range("D19").Value = '43.844.17'
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").Value = '0.0₃0982'
$ws.Range("E20").Value = '  -3.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.82'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.94'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.60'
$ws.Range("E26").Value = '  +28.64%  '
$ws.Range("E27").Value = '  -3.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.90'
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '174.17'
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.91'
$ws.Range("E30").Value = '  +4.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.136'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("E33").Value = '  +3.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.99'
$ws.Range("E34").Value = '  +5.09%  '
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.70'
$ws.Range("E37").Value = '  -3.11%  '
$ws.Range("E38").Value = '  -5.78%  '
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("E40").Value = '  +2.80%  '
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.80'
$ws.Range("E42").Value = '  +4.83%  '
$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.57'
$ws.Range("E43").Value = '  +3.42%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.38'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '98.84'
$ws.Range("E45").Value = '  +1.17%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0951'
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("E48").Value = '  +4.29%  '
$ws.Range("D49").Value = '1.455.60'
$ws.Range("E49").Value = '  -1.47%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.31'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.94'
$ws.Range("E51").Value = '  -4.11%  '
